# Torres_Marlon_ProblemSolving.docx edit:
#   "Moved the header to the left." (commit message) -- in terms of the
#   actual OOXML change, the four header paragraphs ("Marlon Torres",
#   "11/27/2013", "Web Programming Fundamentals - Section 01",
#   "Activity: Problem Solving") go from centre-aligned to right-aligned,
#   and as a side effect of that edit Word's auto-tracked "_GoBack" bookmark
#   (which marks the position of the most recent edit) moves from the
#   empty paragraph right after the header to wrap the header block itself.

$d = $word.ActiveDocument

# 1. Re-justify the four header paragraphs from centered to right-aligned.
#    wdAlignParagraphRight = 2
for ($i = 1; $i -le 4; $i++) {
    $headerPara = $d.Paragraphs.Item($i)
    $headerPara.Range.ParagraphFormat.Alignment = 2
}

# 2. Relocate the hidden "_GoBack" bookmark so it spans the header block
#    (start of paragraph 1 through end of paragraph 4) instead of sitting
#    in the blank paragraph that follows it.
$firstHeaderPara = $d.Paragraphs.Item(1)
$lastHeaderPara = $d.Paragraphs.Item(4)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$goBackRange = $d.Range($firstHeaderPara.Range.Start, $lastHeaderPara.Range.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)
